# Standardize structures of writers amazon:
# Refresh the per-country report rows (file id suffix "_0122_" -> "_0222_",
# r_count, currency label and sum/built_in_total) on Sheet1, rows 2-23.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; A="87811004_0222_RO"; B=447; C="RON"; D=10192.91; E="10192.91"},
    @{Row=3; A="87811004_0222_PE"; B=17; C="PEN"; D=169.05; E="169.05"},
    @{Row=4; A="87811004_0222_HU"; B=1032; C="HUF"; D=2197785; E="2197785"},
    @{Row=5; A="87811004_0222_EU"; B=790; C="EUR"; D=4280.03; E="4280.03"},
    @{Row=6; A="87811004_0222_MX"; B=88; C="MXN"; D=5686.8; E="5686.8"},
    @{Row=7; A="87811004_0222_LL"; B=38; C="USD"; D=100.1; E="100.1"},
    @{Row=8; A="87811004_0222_BG"; B=7; C="BGN"; D=12.18; E="12.18"},
    @{Row=9; A="87811004_0222_BR"; B=41; C="BRL"; D=678.02; E="678.02"},
    @{Row=10; A="87811004_0222_CA"; B=333; C="CAD"; D=1720.6; E="1720.6"},
    @{Row=11; A="87811004_0222_CZ"; B=17; C="CZK"; D=1162; E="1162"},
    @{Row=12; A="87811004_0222_CL"; B=29; C="CLP"; D=43805; E="43805"},
    @{Row=13; A="87811004_0222_CO"; B=24; C="COP"; D=183120; E="183120"},
    @{Row=14; A="87811004_0222_NZ"; B=58; C="NZD"; D=217.42; E="217.42"},
    @{Row=15; A="87811004_0222_AU"; B=357; C="AUD"; D=2286.56; E="2286.56"},
    @{Row=16; A="87811004_0222_CH"; B=67; C="CHF"; D=281.01; E="281.01"},
    @{Row=17; A="87811004_0222_NO"; B=30; C="NOK"; D=887.6; E="887.6"},
    @{Row=18; A="87811004_0222_US"; B=1575; C="USD"; D=11197.9; E="11197.9"},
    @{Row=19; A="87811004_0222_DK"; B=31; C="DKK"; D=478.8; E="478.8"},
    @{Row=20; A="87811004_0222_PL"; B=44; C="PLN"; D=498.12; E="498.12"},
    @{Row=21; A="87811004_0222_SE"; B=54; C="SEK"; D=1740.73; E="1740.73"},
    @{Row=22; A="87811004_0222_JP"; B=32; C="JPY"; D=11340; E="11340"},
    @{Row=23; A="87811004_0222_GB"; B=453; C="GBP"; D=1597.17; E="1597.17"}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A   # file
    $ws.Cells.Item($r, 2).Value = $item.B   # r_count
    $ws.Cells.Item($r, 3).Value = $item.C   # currency
    $ws.Cells.Item($r, 4).Value = $item.D   # sum (number)

    # built_in_total is stored as text even though it mirrors the numeric
    # sum, so force text entry (leading apostrophe) instead of letting the
    # numeric-looking string be auto-converted to a Number.
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value = "'" + $item.E
}
